# Generate Report for Handback
# Adds two new handback entries (aacb56dc-0b34-46ca-b76a-52dcf60615ef and
# f4df2f73-d872-47a0-95be-1856bc9aac1d) as new rows 6 & 7 on every sheet:
#   - "Overview"  (File Name | zh-cn | de-de)
#   - "zh-cn"     (Source File Name | Status | Correspond Handoff File | ...)
#   - "de-de"     (Source File Name | Status | Correspond Handoff File | ...)

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"
$reason = "Include"

$id1 = "aacb56dc-0b34-46ca-b76a-52dcf60615ef"
$id2 = "f4df2f73-d872-47a0-95be-1856bc9aac1d"

$hash1 = "bc46d157e19ad3a880572edd63052e96973858f2"
$hash2 = "29feef96d61737d1a5a4a0597906315098594ffb"

$md1 = "$id1.md"
$md2 = "$id2.md"

$zhXlf1 = "$id1.$hash1.zh-cn.xlf"
$zhXlf2 = "$id2.$hash2.zh-cn.xlf"
$deXlf1 = "$id1.$hash1.de-de.xlf"
$deXlf2 = "$id2.$hash2.de-de.xlf"

$zhHandoffDt1 = "2016-02-26 07:12:53"
$zhHandbackDt1 = "2016-02-26 07:13:57"
$zhHandoffDt2 = "2016-02-26 07:12:53"
$zhHandbackDt2 = "2016-02-26 07:13:57"

$deHandoffDt1 = "2016-02-26 07:13:09"
$deHandbackDt1 = "2016-02-26 07:14:26"
$deHandoffDt2 = "2016-02-26 07:13:09"
$deHandbackDt2 = "2016-02-26 07:14:26"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> rows 6 & 7 : File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/handback/e2e/$md1", "", "", $md1)
$wsOverview.Range("B6").Value = $status
$wsOverview.Range("C6").Value = $status

$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/handback/e2e/$md2", "", "", $md2)
$wsOverview.Range("B7").Value = $status
$wsOverview.Range("C7").Value = $status

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> rows 6 & 7
#   Source File Name | Status | Correspond Handoff File | Correspond Handoff
#   Datetime | Target File | Correspond Handback File | Correspond Handback
#   DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/handback/e2e/$md1", "", "", $md1)
$wsZh.Range("B6").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/handback/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf1", "", "", $zhXlf1)
$wsZh.Range("D6").Value = $zhHandoffDt1
$wsZh.Hyperlinks.Add($wsZh.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/handback/e2e/$md1", "", "", $md1)
$wsZh.Hyperlinks.Add($wsZh.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/handback/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf1", "", "", $zhXlf1)
$wsZh.Range("G6").Value = $zhHandbackDt1
$wsZh.Range("H6").Value = $reason

$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/handback/e2e/$md2", "", "", $md2)
$wsZh.Range("B7").Value = $status
$wsZh.Hyperlinks.Add($wsZh.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/handback/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf2", "", "", $zhXlf2)
$wsZh.Range("D7").Value = $zhHandoffDt2
$wsZh.Hyperlinks.Add($wsZh.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/handback/e2e/$md2", "", "", $md2)
$wsZh.Hyperlinks.Add($wsZh.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/handback/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf2", "", "", $zhXlf2)
$wsZh.Range("G7").Value = $zhHandbackDt2
$wsZh.Range("H7").Value = $reason

# ---------------------------------------------------------------------------
# Sheet "de-de" -> rows 6 & 7 (same layout as zh-cn)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/handback/e2e/$md1", "", "", $md1)
$wsDe.Range("B6").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/handback/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf1", "", "", $deXlf1)
$wsDe.Range("D6").Value = $deHandoffDt1
$wsDe.Hyperlinks.Add($wsDe.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/handback/e2e/$md1", "", "", $md1)
$wsDe.Hyperlinks.Add($wsDe.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/handback/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf1", "", "", $deXlf1)
$wsDe.Range("G6").Value = $deHandbackDt1
$wsDe.Range("H6").Value = $reason

$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/handback/e2e/$md2", "", "", $md2)
$wsDe.Range("B7").Value = $status
$wsDe.Hyperlinks.Add($wsDe.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/handback/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf2", "", "", $deXlf2)
$wsDe.Range("D7").Value = $deHandoffDt2
$wsDe.Hyperlinks.Add($wsDe.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/handback/e2e/$md2", "", "", $md2)
$wsDe.Hyperlinks.Add($wsDe.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/handback/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf2", "", "", $deXlf2)
$wsDe.Range("G7").Value = $deHandbackDt2
$wsDe.Range("H7").Value = $reason

Write-Host "Handback rows added for $id1 and $id2"
